$d = $word.ActiveDocument

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t>The Florida Polytechnic University SGA Policies &amp; Ethics Committee (PEC) is comprised of the PEC Chairman, PEC Vice Chairmen, and all PEC members.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>The Policies &amp; Ethics Committee is responsible for producing, amending and reviewing procedural and statutory bills to ensure they are efficient, clear, and legal.</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>This committee is also responsible for reviewing any grievances filed against an SGA representative and presenting it to the rest of the government.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>Alongside that, we are responsible for listening to student needs and creating legislation to further improve the living and learning experience of those who attend Florida Polytechnic University.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">Contact: Chair John Paul Jones </w:t>
            </w:r>
            <w:r>
              <w:t>jjones5650@floridapoly.com</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Content.InsertXML($xml)
